$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New node row (row 55): "evluator" / metric / evaluate node, mirroring
#    the look of the existing row 54 ("joiner") template directly above it.
# ---------------------------------------------------------------------------

# A55 - bordered on the right, bold-ish header font, wrapped text (same look
# as A54).
$a55 = $ws.Range("A55")
$a55.Value = "evluator"
$a55.Font.Name = "Calibri"
$a55.Font.Size = 10
$a55.Borders.Item(10).LineStyle = 1
$a55.Borders.Item(10).Color = 0
$a55.WrapText = $true

# B55:G55 - boxed cells (left+right thin border), matching B54:G54.
$ws.Range("B55").Value = "metric"
$ws.Range("C55").Value = "evaluate"
$ws.Range("D55").Value = "None"
$ws.Range("E55").Value = '["y_true","y_pred"]'
$ws.Range("F55").Value = '["score"]'
$ws.Range("G55").Value = "evaluate"

foreach ($col in @("B", "C", "D", "E", "F", "G")) {
    $cell = $ws.Range($col + "55")
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(7).Color = 0
    $cell.Borders.Item(10).LineStyle = 1
    $cell.Borders.Item(10).Color = 0
}

# ---------------------------------------------------------------------------
# 2. Re-stamp row 54's own formatting (format-only self copy) - this is the
#    bit of tidy-up Excel performed on the "joiner" row while the new row was
#    being authored, normalising its cell formats.
# ---------------------------------------------------------------------------
$ws.Range("A54:G54").Copy() | Out-Null
$ws.Range("A54:G54").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Update the sheet view so the just-entered row is in frame, matching the
#    recorded cursor / scroll position after authoring the new row.
# ---------------------------------------------------------------------------
$ws.Range("C58").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 44
$win.ScrollColumn = 2
